# Updated cryptos list on Sat Sep 14 10:38:30 UTC 2024 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns of the cryptos table
# in-place, keeping every touched cell a plain text value (matching the
# existing inline-string cells) rather than letting Excel re-interpret
# numeric-looking price strings (e.g. "6.20") as numbers, which would
# silently drop significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.971.02"
$ws.Range("E2").Value = "  +3.06%  "

$ws.Range("D3").Value = "2.418.37"
$ws.Range("E3").Value = "  +2.47%  "

$ws.Range("E4").Value = "  +0.05%  "

# D5 reads like a number ("551.48"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D5").Value = "'551.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "

# D6 reads like a number ("137.22"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D6").Value = "'137.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.62%  "

$ws.Range("E7").Value = "  +0.02%  "

# D8 reads like a number ("0.583"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D8").Value = "'0.583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.98%  "

$ws.Range("E9").Value = "  -0.32%  "

$ws.Range("E10").Value = "  +3.92%  "

$ws.Range("E11").Value = "  -1.72%  "

# D12 reads like a number ("0.357"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D12").Value = "'0.357"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "

# D13 reads like a number ("24.77"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D13").Value = "'24.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.94%  "

$ws.Range("D14").Value = "2.848.65"
$ws.Range("E14").Value = "  +2.55%  "

$ws.Range("D15").Value = "59.918.32"
$ws.Range("E15").Value = "  +3.08%  "

# D16 reads like a number ("0.0000137"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D16").Value = "'0.0000137"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").Value = "2.428.98"
$ws.Range("E17").Value = "  +3.51%  "

# D18 reads like a number ("11.29"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D18").Value = "'11.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.20%  "

$ws.Range("E19").Value = "  +1.34%  "

# D20 reads like a number ("331.07"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D20").Value = "'331.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("E21").Value = "  -2.10%  "

$ws.Range("E22").Value = "  -0.09%  "

# D23 reads like a number ("65.87"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D23").Value = "'65.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.70%  "

# D24 reads like a number ("0.173"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D24").Value = "'0.173"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.79%  "

# D25 reads like a number ("8.59"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D25").Value = "'8.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.12%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").Value = "  +1.48%  "

$ws.Range("E28").Value = "  +5.63%  "

$ws.Range("E29").Value = "  +1.21%  "

# D30 reads like a number ("170.58"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D30").Value = "'170.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.30%  "

# D31 reads like a number ("6.20"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D31").Value = "'6.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.99%  "

$ws.Range("E32").Value = "  +1.33%  "

$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("E35").Value = "  +4.35%  "

$ws.Range("E36").Value = "  +0.09%  "

# D37 reads like a number ("4.18"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D37").Value = "'4.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("E38").Value = "  +0.27%  "

# D39 reads like a number ("39.49"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D39").Value = "'39.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.95%  "

$ws.Range("E40").Value = "  +0.74%  "

# D41 reads like a number ("314.49"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D41").Value = "'314.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.97%  "

$ws.Range("E42").Value = "  +0.18%  "

# D43 reads like a number ("138.52"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D43").Value = "'138.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.54%  "

# D44 reads like a number ("0.0962"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D44").Value = "'0.0962"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.26%  "

# D45 reads like a number ("0.0519"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D45").Value = "'0.0519"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("E46").Value = "  +2.30%  "

# D47 reads like a number ("19.30"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D47").Value = "'19.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.68%  "

# D48 reads like a number ("0.410"); force a leading quote so Excel
# keeps it as text, then reset the style so no explicit text format lingers
$ws.Range("D48").Value = "'0.410"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.79%  "

$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("E50").Value = "  +0.79%  "

$ws.Range("E51").Value = "  -0.44%  "
